$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.078.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.25%  '
$ws.Range("D3").Value = "'3.454.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.28%  '
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = "'581.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").Value = "'165.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.23%  '
$ws.Range("D7").Value = "'0.602"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.02%  '
$ws.Range("D8").Value = "'3.441.03"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.45%  '
$ws.Range("D10").Value = "'0.185"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.94%  '
$ws.Range("D11").Value = "'6.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.08%  '
$ws.Range("D12").Value = "'0.568"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.11%  '
$ws.Range("D13").Value = "'45.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.54%  '
$ws.Range("D14").Value = "'0.0000270"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.67%  '
$ws.Range("D15").Value = "'3.989.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.87%  '
$ws.Range("D16").Value = "'611.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -10.38%  '
$ws.Range("D17").Value = "'8.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -7.84%  '
$ws.Range("D18").Value = "'3.449.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.78%  '
$ws.Range("D19").Value = "'67.889.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.62%  '
$ws.Range("E20").Value = '  -3.42%  '
$ws.Range("D21").Value = "'17.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.19%  '
$ws.Range("D22").Value = "'10.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.15%  '
$ws.Range("D23").Value = "'0.867"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.72%  '
$ws.Range("D24").Value = "'15.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.99%  '
$ws.Range("D25").Value = "'95.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.24%  '
$ws.Range("D26").Value = "'3.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.57%  '
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").Value = "'2.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.80%  '
$ws.Range("D29").Value = "'8.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.97%  '
$ws.Range("D30").Value = "'32.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.30%  '
$ws.Range("D31").Value = "'8.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.18%  '
$ws.Range("E32").Value = '  -6.67%  '
$ws.Range("D33").Value = "'1.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.61%  '
$ws.Range("D34").Value = "'6.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -9.42%  '
$ws.Range("D35").Value = "'582.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.93%  '
$ws.Range("D36").Value = "'10.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.78%  '
$ws.Range("D37").Value = "'56.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.88%  '
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").Value = "'0.1000"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.22%  '
$ws.Range("D40").Value = "'3.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -14.87%  '
$ws.Range("D41").Value = "'0.135"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.33%  '
$ws.Range("D42").Value = "'0.0429"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.99%  '
$ws.Range("D43").Value = "'3.355.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.70%  '
$ws.Range("D44").Value = "'0.318"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.47%  '
$ws.Range("D45").Value = "'32.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.13%  '
$ws.Range("D46").Value = "'0.0₃0677"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.15%  '
$ws.Range("D47").Value = "'2.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.08%  '
$ws.Range("D48").Value = "'2.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.98%  '
$ws.Range("E49").Value = '  -5.14%  '
$ws.Range("D50").Value = "'132.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.00%  '
$ws.Range("D51").Value = "'5.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.47%  '
